# Update "想去人数" (interest count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 506
$wsExpo.Range("F3").Value = 5990
$wsExpo.Range("F6").Value = 111
$wsExpo.Range("F9").Value = 553

# 演出 (Performances) sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 5

# 全部类型 (All types, combined) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 506
$wsAll.Range("F3").Value = 5990
$wsAll.Range("F5").Value = 5
$wsAll.Range("F7").Value = 111
$wsAll.Range("F11").Value = 553
